$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "bibi" base numbers for the last row (ano_obj 2025):
# total_customers (C8) and new_customers (E8) increased; returning_customers (D8) unchanged.
$ws.Range("C8").Value = 1330
$ws.Range("E8").Value = 1120

# Recompute dependent rate columns (G = retention_rate, H = returning_rate)
# from the refreshed totals; new_rate (F8) is unaffected by this update.
$c8 = $ws.Range("C8").Value()
$d8 = $ws.Range("D8").Value()
$e8 = $ws.Range("E8").Value()

$ws.Range("G8").Value = $e8 / $c8 * 100
$ws.Range("H8").Value = $d8 / $c8 * 100
